# Add a new quarterly column (BB) to the right of the existing data (which
# currently ends at column BA), mirroring how prior quarter columns were
# appended: the header row gets a new date serial, and every data row's
# new cell is populated with the same value as its current last (BA) cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header date (row 1, column BB = column 54). Copy the previous
# header cell (BA1) so the new cell picks up the same style/number
# format, then overwrite with the new quarter's date serial.
$ws.Cells.Item(1, 53).Copy($ws.Cells.Item(1, 54))
$ws.Cells.Item(1, 54).Value2 = 45986

# Copy each data row's existing last value (column BA = 53) into the new
# column BB (54) for rows that have data there, preserving formatting.
for ($r = 3; $r -le 21; $r++) {
    $srcCell = $ws.Cells.Item($r, 53)
    if ($srcCell.Value2 -ne $null) {
        $srcCell.Copy($ws.Cells.Item($r, 54))
    }
}
